# 29 Mayıs verileri eklendi.
# Append the 2020-05-28 (serial 43979) COVID-19 daily stats row to the
# "data" worksheet and extend the Table3 Excel table to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data row to append right after the existing last row (row 77).
$newRow = 78
$ws.Range("A$newRow").Value = 43979
$ws.Range("B$newRow").Value = 33559
$ws.Range("C$newRow").Value = 1182
$ws.Range("D$newRow").Value = 30
$ws.Range("E$newRow").Value = 1576

# Grow the worksheet's Excel Table ("Table3") so it covers the new row too;
# this keeps the table's ref/autoFilter in sync with the appended data.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E$newRow"))

# Match the workbook's saved selection state after the edit.
$ws.Range("E77").Select() | Out-Null
